$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.123.03"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").Value = "1.656.72"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.96"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5261"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.03%  "
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2606"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06349"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.45"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07783"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.508"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.12%  "
$ws.Range("D13").Value = "1.658.00"
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5488"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.97%  "
$ws.Range("D15").Value = "0.0₅8233"
$ws.Range("E15").Value = "  +1.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.39"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.47%  "
$ws.Range("D17").Value = "26.138.07"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("E18").Value = "  -0.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.579"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "192.17"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.09"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "141.71"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1250"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.277"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.19"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.435"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05918"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.279"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.32%  "
$ws.Range("E31").Value = "  -1.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.260"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.585"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9539"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.790"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.47%  "
$ws.Range("E36").Value = "  -0.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5705"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01619"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.02%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8490"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.80%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.782"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.002"
$ws.Range("D41").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "103.06"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.93%  "
$ws.Range("D43").Value = "1.030.39"
$ws.Range("E43").Value = "  +1.80%  "
$ws.Range("D44").Value = "1.802.47"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "57.46"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.004"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.53%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4301"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.04%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.483"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.97%  "
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("E50").Value = "  -1.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.09728"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.43%  "
